# Updates crypto price/volume figures in the "cryptos" list on Sheet1.
# For D-column values that look like plain numbers (e.g. "602.41"), we
# briefly force a text number format so Excel keeps them as literal
# strings (matching the source data, which stores prices as text because
# some of them use "." as a thousands separator, e.g. "70.298.99") and
# then restore the default "Normal" style so no stray formatting is left
# behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.298.99"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").Value = "3.617.50"
$ws.Range("E3").Value = "  +2.15%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "195.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.34%  "

$ws.Range("E7").Value = "  -0.75%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("E9").Value = "  +3.52%  "

$ws.Range("E10").Value = "  -0.95%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.28"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000305"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.88%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.19%  "

$ws.Range("D14").Value = "4.193.27"
$ws.Range("E14").Value = "  +2.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "602.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.70%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "12.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.39%  "

$ws.Range("D17").Value = "70.412.63"
$ws.Range("E17").Value = "  +0.42%  "

$ws.Range("D18").Value = "3.636.46"
$ws.Range("E18").Value = "  +3.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.80%  "

$ws.Range("E20").Value = "  +1.70%  "

$ws.Range("E21").Value = "  +0.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.75%  "

$ws.Range("E23").Value = "  -1.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "102.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("E30").Value = "  +9.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.57%  "

$ws.Range("E32").Value = "  -2.73%  "

$ws.Range("E33").Value = "  +1.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.26"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("D35").Value = "0.0₃0880"
$ws.Range("E35").Value = "  +2.37%  "

$ws.Range("D36").Value = "3.934.16"
$ws.Range("E36").Value = "  +5.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "534.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.11%  "

$ws.Range("E38").Value = "  +0.12%  "

$ws.Range("E39").Value = "  +0.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.36%  "

$ws.Range("E41").Value = "  -1.27%  "

$ws.Range("E42").Value = "  -2.41%  "

$ws.Range("E43").Value = "  +0.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0462"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.26%  "

$ws.Range("E47").Value = "  -0.55%  "

$ws.Range("E48").Value = "  -0.55%  "

$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000250"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.52%  "

$ws.Range("E51").Value = "  +1.25%  "
